$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells (prices, volume %, and shifted coin rows) per diff
$ws.Cells.Item(2, 4).Value = '29.467.72'
$ws.Cells.Item(2, 5).Value = '  +0.41%  '
$ws.Cells.Item(3, 4).Value = '1.852.47'
$ws.Cells.Item(3, 5).Value = '  +0.50%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '240.90'
$ws.Cells.Item(5, 5).Value = '  +0.95%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.6305'
$ws.Cells.Item(6, 5).Value = '  +0.57%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.002'
$ws.Cells.Item(7, 5).Value = '  +0.11%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.07662'
$ws.Cells.Item(8, 5).Value = '  +1.66%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.2940'
$ws.Cells.Item(9, 5).Value = '  +0.02%  '
$ws.Cells.Item(10, 5).Value = '  +0.52%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07759'
$ws.Cells.Item(11, 5).Value = '  +0.86%  '
$ws.Cells.Item(12, 4).Value = '1.860.86'
$ws.Cells.Item(12, 5).Value = '  +0.63%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.031'
$ws.Cells.Item(13, 5).Value = '  +1.21%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.6808'
$ws.Cells.Item(14, 5).Value = '  +0.58%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.00001065'
$ws.Cells.Item(15, 5).Value = '  +4.40%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '83.65'
$ws.Cells.Item(16, 5).Value = '  +0.82%  '
$ws.Cells.Item(17, 4).Value = '2.113.37'
$ws.Cells.Item(17, 5).Value = '  +0.55%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '6.171'
$ws.Cells.Item(18, 5).Value = '  +0.95%  '
$ws.Cells.Item(19, 4).Value = '29.482.28'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '229.13'
$ws.Cells.Item(20, 5).Value = '  +0.62%  '
$ws.Cells.Item(21, 5).Value = '  +0.53%  '
$ws.Cells.Item(22, 5).Value = '  +0.08%  '
$ws.Cells.Item(23, 5).Value = '  -0.43%  '
$ws.Cells.Item(24, 5).Value = '  +0.10%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '156.96'
$ws.Cells.Item(25, 5).Value = '  +0.11%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.1384'
$ws.Cells.Item(26, 5).Value = '  -0.30%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '8.401'
$ws.Cells.Item(27, 5).Value = '  +0.71%  '
$ws.Cells.Item(28, 5).Value = '  +0.65%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.325'
$ws.Cells.Item(29, 5).Value = '  +4.85%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.471'
$ws.Cells.Item(30, 5).Value = '  +0.97%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.05677'
$ws.Cells.Item(31, 5).Value = '  +1.14%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.134'
$ws.Cells.Item(32, 5).Value = '  +0.46%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.046'
$ws.Cells.Item(33, 5).Value = '  +0.41%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.849'
$ws.Cells.Item(34, 5).Value = '  +0.94%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.165'
$ws.Cells.Item(35, 5).Value = '  +1.06%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.7039'
$ws.Cells.Item(36, 5).Value = '  -1.05%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.588'
$ws.Cells.Item(37, 5).Value = '  -0.22%  '
$ws.Cells.Item(38, 5).Value = '  +0.52%  '
$ws.Cells.Item(39, 5).Value = '  -0.47%  '
$ws.Cells.Item(40, 4).Value = '1.219.10'
$ws.Cells.Item(40, 5).Value = '  -1.96%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '6.555'
$ws.Cells.Item(41, 5).Value = '  +5.92%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.9107'
$ws.Cells.Item(42, 5).Value = '  +0.95%  '
$ws.Cells.Item(43, 5).Value = '  +0.18%  '
$ws.Cells.Item(44, 2).Value = 'Quant'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '101.71'
$ws.Cells.Item(44, 5).Value = '  +0.14%  '
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '66.36'
$ws.Cells.Item(45, 5).Value = '  +0.41%  '
$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.00000000121'
$ws.Cells.Item(46, 5).Value = '  -0.87%  '
$ws.Cells.Item(47, 2).Value = 'Aptos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.111'
$ws.Cells.Item(47, 5).Value = '  +0.33%  '
$ws.Cells.Item(48, 2).Value = 'TheSandbox'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.4021'
$ws.Cells.Item(48, 5).Value = '  +0.81%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '9.027'
$ws.Cells.Item(49, 5).Value = '  +0.50%  '
$ws.Cells.Item(50, 2).Value = 'RenderToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.682'
$ws.Cells.Item(50, 5).Value = '  +0.13%  '
$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.1138'
$ws.Cells.Item(51, 5).Value = '  +2.12%  '
